$wb = $excel.ActiveWorkbook

# --- Sheet 1: "总计" (totals) ---------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# The existing 2022-Q2 row slides down to row 3 (same values it already
# had), and the new 2022-Q3 figures take over row 2 (mirroring "a new
# quarter was inserted above the previous one").
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 0.29

# Give the new A3 cell the same look as A2 (bold/centered header style).
$totals.Range("A2").Copy() | Out-Null
$totals.Range("A3").PasteSpecial(-4122) | Out-Null

$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 0.23

# --- New sheet: "2022-Q3" --------------------------------------------------
# Inserted right after "总计", ahead of the existing "2022-Q2" sheet. Clone
# "总计" (rather than starting from a blank sheet) so the new tab picks up
# the same sheet setup / page margins / header style, then wipe the cloned
# values back out before filling in the fund table.
$totals.Copy($null, $totals)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"
$q3.Cells.ClearContents()

# Extend the header styling (copied from "总计"'s B1:D1/A2) across the
# wider fund table - columns E1:H1 and the second index cell A3.
$q3.Range("B1:D1").Copy() | Out-Null
$q3.Range("E1:H1").PasteSpecial(-4122) | Out-Null
$q3.Range("A2").Copy() | Out-Null
$q3.Range("A3").PasteSpecial(-4122) | Out-Null

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("C2").Value = "嘉实北交所精选两年定期混合A"
$q3.Range("H2").Value = 1

$q3.Range("A3").Value = 1
$q3.Range("C3").Value = "嘉实北交所精选两年定期混合C"
$q3.Range("H3").Value = 1

# These columns hold digit-strings (fund code / percentages) that must stay
# text (e.g. "014269" keeps its leading zero) rather than becoming numbers,
# so format the cells as text before typing the values in - exactly what
# Excel needs to keep a quoted-looking number as a string.
# (Multi-area "A,B" ranges only apply formatting to their first area here,
# so the text columns are handled as separate contiguous ranges instead.)
$q3.Range("B2:B3").NumberFormat = "@"
$q3.Range("D2:G3").NumberFormat = "@"

$q3.Range("B2").Value = "014269"
$q3.Range("D2").Value = "2.72"
$q3.Range("E2").Value = "90.37"
$q3.Range("F2").Value = "7.21"
$q3.Range("G2").Value = "0.1961"

$q3.Range("B3").Value = "014270"
$q3.Range("D3").Value = "0.53"
$q3.Range("E3").Value = "90.37"
$q3.Range("F3").Value = "7.21"
$q3.Range("G3").Value = "0.0382"

# Drop back to the default style now that the values are locked in as text
# (keeps the cells' look/format identical to a plain, unstyled cell).
$q3.Range("B2:B3").Style = "Normal"
$q3.Range("D2:G3").Style = "Normal"
